# Handback status report generation:
# Remove the second data row (the "48bc4e3f-4ec2-44de-a3ab-d3ced75e777a" entry)
# from every worksheet, and refresh the handback timestamps for the
# remaining ("01a7335d-4021-43b2-ba3d-305a29b44724") entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop the hyperlinks that live on the row we are about to remove, then
# delete the row itself (this shifts nothing else, since it's the last row).
$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(3).Delete()

# Restore the hyperlink that belongs to the row which remains.
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/080e74cdff42a93e842e66ce59520e7b4567de6b/e2e/01a7335d-4021-43b2-ba3d-305a29b44724.md", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/080e74cdff42a93e842e66ce59520e7b4567de6b/e2e/01a7335d-4021-43b2-ba3d-305a29b44724.md", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a67fe45125b6b83f447657c1f4beceb8c59aa393/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.zh-cn.xlf", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/abf348f3ed26f633f3cccd3aa4455a6926cfb973/e2e/01a7335d-4021-43b2-ba3d-305a29b44724.md", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7b469738e2d27a4241814041f3ce6d92e971df46/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.zh-cn.xlf", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.zh-cn.xlf") | Out-Null

# Refresh the handoff / handback timestamps for the remaining row.
$ws2.Range("E2").Value = "2016-03-20 04:47:19"
$ws2.Range("H2").Value = "2016-03-20 04:48:02"

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/080e74cdff42a93e842e66ce59520e7b4567de6b/e2e/01a7335d-4021-43b2-ba3d-305a29b44724.md", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f767765f93a07b8a747de7089792b0b648e4e38/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.de-de.xlf", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/493bc7398acaa3038baa524b4cb6648133748f3b/e2e/01a7335d-4021-43b2-ba3d-305a29b44724.md", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/56d75d80323a0399677aa71921884893a148d9fe/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.de-de.xlf", "", "", "01a7335d-4021-43b2-ba3d-305a29b44724.f7c95888a78221389ccea09f2aec87ee623b8f85.de-de.xlf") | Out-Null

$ws3.Range("E2").Value = "2016-03-20 04:47:27"
$ws3.Range("H2").Value = "2016-03-20 04:48:15"

$wb.Save()
